$d = $word.ActiveDocument

# The document has Track Changes enabled; disable it so our edits apply
# directly instead of being recorded as insertions/deletions.
$d.TrackRevisions = $false

# --- Delete all existing comments (id 0, 1, 3) -----------------------------
# Comment 0: "Hora esta mal. Ve o mail do prof" (wraps the meeting time range)
# Comment 1: "Devias planear uma reuniao de 60min" (on "Agenda:")
# Comment 3: "Data mal." (wraps the next-meeting date)
while ($d.Comments.Count -gt 0) {
    $d.Comments(1).Delete()
}

# --- Fix the meeting time: 22:00-22:45 -> 21:40-22:40 ----------------------
$null = $d.Paragraphs(1).Range.Find.Execute(
    "22:00-22:45", $true, $false, $false, $false, $false, $true, 1, $false,
    "21:40-22:40", 2)

# --- Swap roles: Time Keeper / Note Keeper ----------------------------------
$null = $d.Paragraphs(10).Range.Find.Execute(
    "Rui Ganhoto", $true, $false, $false, $false, $false, $true, 1, $false,
    "João Girão", 2)

$null = $d.Paragraphs(11).Range.Find.Execute(
    "João Girão", $true, $false, $false, $false, $false, $true, 1, $false,
    "Mário Oliveira", 2)

# --- Agenda timing tweaks ----------------------------------------------------
# "Analyze use cases and mockups[5 min];" -> "...[10 min];"
$null = $d.Paragraphs(16).Range.Find.Execute(
    "5 min", $true, $false, $false, $false, $false, $true, 1, $false,
    "10 min", 2)

# "Analyze requirements [10 min];" -> "...[15 min];"
$null = $d.Paragraphs(17).Range.Find.Execute(
    "10 min", $true, $false, $false, $false, $false, $true, 1, $false,
    "15 min", 2)

# "Goals for next week and any other business [5 min];" -> "...[10 min];"
$null = $d.Paragraphs(20).Range.Find.Execute(
    "5 min", $true, $false, $false, $false, $false, $true, 1, $false,
    "10 min", 2)

# --- Fix the Next Meeting date: 22/04/2013 -> 29/04/2013 -------------------
$null = $d.Paragraphs(26).Range.Find.Execute(
    "22/04/2013", $true, $false, $false, $false, $false, $true, 1, $false,
    "29/04/2013", 2)
